$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Structural change: split the old "Size" column into two columns
# ("Rated Power" / "Rated Energy"), and add a new "Average Annual
# Resting SOC" column before "Site Ambient Temperature Range".
# ------------------------------------------------------------------

# Insert a new column at C (old C "Size" shifts right to D, becoming
# "Rated Energy"); and a new column at L (old K "Site Ambient..."
# shifts right to M).
$ws.Columns("C").Insert()
$ws.Columns("L").Insert()

# New columns pick up a manual width close to their left neighbour's
# (mirrors what Excel does cosmetically when a column is inserted).
$ws.Columns("C").ColumnWidth = 26
$ws.Columns("L").ColumnWidth = 33.8

# ------------------------------------------------------------------
# Header row
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Rated Power"
$ws.Range("D1").Value = "Rated Energy"
$ws.Range("L1").Value = "Average Annual Resting SOC"

# ------------------------------------------------------------------
# Row 2 - Golden Triangle II: "50 MW / 200 MWh" -> "50 MW" + "200 MWh"
# ------------------------------------------------------------------
$ws.Range("C2").Value = "50 MW"
$ws.Range("D2").Value = "200 MWh"

# ------------------------------------------------------------------
# Row 3 - Happy Valley: "82 MW / 328 MWh" -> "82 MW" + "328 MWh"
# ------------------------------------------------------------------
$ws.Range("C3").Value = "82 MW"
$ws.Range("D3").Value = "328 MWh"

# Row 3's resting-SOC cell: plain (single-run) text. Written first so
# it lands in its own shared-string slot, distinct from row 2's
# rich-text version below (same visible text, different formatting).
$ws.Range("L3").Value = "≤50%"
$ws.Range("L3").NumberFormat = "0%"

# Row 2's resting-SOC cell: rich text "<=50%" (the "<=" glyph renders in
# a fallback font, "Aptos Narrow", the rest stays in Calibri) formatted
# with a percent number format.
$ws.Range("L2").Value = "≤50%"
$ws.Range("L2").NumberFormat = "0%"
$ws.Range("L2").Characters(1, 1).Font.Name = "Aptos Narrow"
$ws.Range("L2").Characters(2, 3).Font.Name = "Calibri"

# ------------------------------------------------------------------
# Cosmetic touch-ups to roughly track the authored view state
# (scroll so column E is left-most visible, select L8).
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L8").Select() | Out-Null
